# Update Enfermeria worksheet: roll the "Mes" column from 2023-06-01 to
# 2023-07-01 and refresh the "Cantidad" figures for the July report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New month value for column B (all data rows): 2023-07-01 == serial 45108
$newDate = 45108

# row => new Cantidad (D) value; rows not listed keep their prior value
$quantities = @{
    2  = 100
    3  = 136
    4  = 10
    5  = 2
    6  = 5
    7  = 15
    8  = 2
    9  = 60
    10 = 0
    12 = 122
    13 = 0
    14 = 0
    15 = 0
    16 = 64
    17 = 66
    18 = 61
    19 = 3
    20 = 1
}

for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 2).Value = $newDate
    if ($quantities.ContainsKey($row)) {
        $ws.Cells.Item($row, 4).Value = $quantities[$row]
    }
}

# Move the saved selection from D14 to D13, matching the sheetView state.
$ws.Range("D13").Select()
